$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new daily price record was inserted above the existing row 97 ("Puerro"
# price list for Vega Modelo de Temuco), shifting all subsequent rows (old
# 97..224) down by one (to 98..225). Insert a fresh row at 97 and populate it
# with the new record's values.
$ws.Rows.Item(97).Insert()

$ws.Cells.Item(97, 1).Value = 10
$ws.Cells.Item(97, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(97, 3).Value = "La Araucanía"
$ws.Cells.Item(97, 4).Value = 44781
$ws.Cells.Item(97, 5).Value = 9
$ws.Cells.Item(97, 6).Value = 100112005
$ws.Cells.Item(97, 7).Value = "Puerro"
$ws.Cells.Item(97, 8).Value = "Azul de Maquehue"
$ws.Cells.Item(97, 9).Value = "Primera"
$ws.Cells.Item(97, 10).Value = 40
$ws.Cells.Item(97, 11).Value = 16000
$ws.Cells.Item(97, 12).Value = 16000
$ws.Cells.Item(97, 13).Value = 16000
$ws.Cells.Item(97, 14).Value = "$/docena de paquetes"
$ws.Cells.Item(97, 15).Value = "Provincia de Cautín"
$ws.Cells.Item(97, 16).Value = 1333
$ws.Cells.Item(97, 17).Value = 12
$ws.Cells.Item(97, 18).Value = "Hortaliza"
